$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 82, shifting the existing rows 82:93 down to 84:95.
$ws.Range("A82:T83").Insert()

# New row 82: Chirimoya "Especial", Provincia del Elquí, fecha 2021-11-05 (serial 44505)
$ws.Range("A82").Value = 10
$ws.Range("B82").Value = "Vega Modelo de Temuco"
$ws.Range("C82").Value = "La Araucanía"
$ws.Range("D82").Value = 44505
$ws.Range("E82").Value = 9
$ws.Range("F82").Value = "Fruta"
$ws.Range("G82").Value = 100107
$ws.Range("H82").Value = "Otros"
$ws.Range("I82").Value = 100107002
$ws.Range("J82").Value = "Chirimoya"
$ws.Range("K82").Value = "Cultivar IV Región"
$ws.Range("L82").Value = "Especial"
$ws.Range("M82").Value = 20
$ws.Range("N82").Value = 3500
$ws.Range("O82").Value = 3500
$ws.Range("P82").Value = 3500
$ws.Range("Q82").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R82").Value = "Provincia del Elquí"
$ws.Range("S82").Value = 3500
$ws.Range("T82").Value = 1

# New row 83: Chirimoya "Primera", Provincia del Elquí, fecha 2021-11-05 (serial 44505)
$ws.Range("A83").Value = 10
$ws.Range("B83").Value = "Vega Modelo de Temuco"
$ws.Range("C83").Value = "La Araucanía"
$ws.Range("D83").Value = 44505
$ws.Range("E83").Value = 9
$ws.Range("F83").Value = "Fruta"
$ws.Range("G83").Value = 100107
$ws.Range("H83").Value = "Otros"
$ws.Range("I83").Value = 100107002
$ws.Range("J83").Value = "Chirimoya"
$ws.Range("K83").Value = "Cultivar IV Región"
$ws.Range("L83").Value = "Primera"
$ws.Range("M83").Value = 35
$ws.Range("N83").Value = 3000
$ws.Range("O83").Value = 3000
$ws.Range("P83").Value = 3000
$ws.Range("Q83").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R83").Value = "Provincia del Elquí"
$ws.Range("S83").Value = 3000
$ws.Range("T83").Value = 1
